# Updated symbol list on Wed Jan  4 09:37:55 UTC 2023 with GitHub Actions
#
# Refreshes the "cryptos" price table: per-row Price (col D) and
# Volume(1h) (col E) values, plus a few rows (9-14) whose Coin/Link/Price
# pairs shifted because the underlying ranking reshuffled.
#
# Price/Volume cells are stored as plain text in the workbook (e.g. "255.23",
# "3.71%"), not as numbers - so numeric-looking values are written with a
# leading apostrophe to force Excel to keep them as text instead of
# re-interpreting them as a Number/Percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'255.23"
$ws.Range("E2").Value = "'3.71%"

# Row 3 - OKB
$ws.Range("D3").Value = "'28.14"
$ws.Range("E3").Value = "'-5.31%"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'5.248"
$ws.Range("E4").Value = "'1.79%"

# Row 5 - Cronos
$ws.Range("E5").Value = "'1.45%"

# Row 6 - KuCoinToken
$ws.Range("D6").Value = "'6.706"
$ws.Range("E6").Value = "'0.76%"

# Row 7 - MXToken
$ws.Range("D7").Value = "'0.8682"
$ws.Range("E7").Value = "'2.27%"

# Row 8 - FTXToken
$ws.Range("D8").Value = "'1.035"
$ws.Range("E8").Value = "'21.10%"

# Row 9 - was "One", now "WazirX"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1410"
$ws.Range("E9").Value = "'1.80%"

# Row 10 - was "WazirX", now "MandalaExchangeToken"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07115"
$ws.Range("E10").Value = "'0.43%"

# Row 11 - was "MandalaExchangeToken", now "BitrueCoin"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.03188"
$ws.Range("E11").Value = "'-1.97%"

# Row 12 - was "BitrueCoin", now "BitMartToken"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "'0.09226"
$ws.Range("E12").Value = "'-1.61%"

# Row 13 - was "BitMartToken", now "BitForexToken"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "'0.001544"
$ws.Range("E13").Value = "'1.27%"

# Row 14 - was "BitForexToken", now "One"
$ws.Range("B14").Value = "One"
$ws.Range("C14").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D14").Value = "'0.0006042"
$ws.Range("E14").Value = "'0.90%"

# Row 15 - TigerCash
$ws.Range("D15").Value = "'0.005830"
$ws.Range("E15").Value = "'-0.86%"

# Row 17 - GateToken
$ws.Range("E17").Value = "'-0.26%"

# Row 18 - BTSEToken
$ws.Range("D18").Value = "'2.223"
$ws.Range("E18").Value = "'1.63%"

# Row 19 - BitpandaEcosystemToken
$ws.Range("E19").Value = "'0.50%"

# Row 20 - LiechtensteinCryptoassetsExchange
$ws.Range("D20").Value = "'0.03476"
$ws.Range("E20").Value = "'3.53%"

# Row 22 - MCDex
$ws.Range("D22").Value = "'3.549"
$ws.Range("E22").Value = "'1.77%"

# Row 23 - CoinExToken
$ws.Range("D23").Value = "'0.04167"
$ws.Range("E23").Value = "'0.84%"

# Row 24 - ZBToken
$ws.Range("E24").Value = "'-4.42%"

# Row 25 - BitKan
$ws.Range("D25").Value = "'0.001223"
$ws.Range("E25").Value = "'-0.13%"

# Row 26 - HotbitToken
$ws.Range("D26").Value = "'0.004881"
$ws.Range("E26").Value = "'17.86%"

# Row 27
$ws.Range("E27").Value = "'0.10%"

# Row 28
$ws.Range("D28").Value = "'0.0001939"
$ws.Range("E28").Value = "'33.83%"

# Row 40 - IDEX
$ws.Range("D40").Value = "'0.03818"
$ws.Range("E40").Value = "'1.85%"

# Row 41 - KickToken
$ws.Range("D41").Value = "'0.005732"
$ws.Range("E41").Value = "'0.99%"

# Row 42 - BKEXToken
$ws.Range("D42").Value = "'0.1103"
$ws.Range("E42").Value = "'3.05%"

# Row 43 - CEJI
$ws.Range("D43").Value = "'0.002339"
$ws.Range("E43").Value = "'1.76%"

# Row 44 - LocalTraders
$ws.Range("D44").Value = "'0.009681"
$ws.Range("E44").Value = "'9.04%"

# Row 45 - CoinLion
$ws.Range("D45").Value = "'0.00005235"
$ws.Range("E45").Value = "'-4.84%"

# Row 46 - Kangarootoken
$ws.Range("E46").Value = "'0.10%"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = "'0.09303"
$ws.Range("E47").Value = "'31.06%"

# Row 48 - BOLO
$ws.Range("E48").Value = "'-12.80%"

# Row 49 - CryptobidCoin
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.10%"

# Row 50 - SpecialPowerGold
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.10%"
